# Regenerate s_val data to filter save games.
# This updates the computed TB/d2S/K/IP/sum columns (B:E, G) for rows 2-8.
# The "Win" column (F) is unaffected by this regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    3 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068,  21.98653043760045)
    4 = @(1.445647641019636, 9.983522426115931, 0.7210945179870265, 13.86384647080068,  26.01411105592328)
    5 = @(0.04172184405617529, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.337238039619874)
    6 = @(3.272327238179451, 9.983522426115931, 0.7210945179870265, 13.86384647080068,  27.84079065308309)
    7 = @(0.04172184405617529, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.351702369198972)
    8 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G: sum
}
